$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), matching style of existing header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new column I and J values for rows 2-11
$values = @{
    2  = @(1, 6)
    3  = @(1, 6)
    4  = @(1, 6)
    5  = @(1, 5)
    6  = @(2, 6)
    7  = @(4, 8)
    8  = @(1, 5)
    9  = @(1, 5)
    10 = @(9, 9)
    11 = @(7, 7)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
